$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell H1, matching style of G1
$ws.Range("G1").Copy() | Out-Null
$ws.Range("H1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("H1").Value = "Save"

# New data cells
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
